# Mark additional functions on the "Functions" sheet as implemented
# ("N/A" -> "Done" in the Status column), and move the sheet's
# selection/scroll to D7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functions")

# Rows whose Status (column D) flips from "N/A" to "Done".
$rows = @(6, 8, 9, 11, 12, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = "Done"
}

# Update the visible selection / scroll position for the sheet.
$ws.Activate()
$ws.Range("D7").Select()
